# cap nhat ho so simple
#
# The original row 4 entry (vehicle plate "84A-142.58" that departed/returned
# on 2025-11-18) is replaced by a new entry for plate "84A-245.53" that
# departed 2025-11-25 07:00:00 and returned 2025-11-25 14:00:00, logging 7
# total hours.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B4").Value = "84A-245.53"
$ws.Range("C4").Value = "2025-11-25 07:00:00"
$ws.Range("D4").Value = "2025-11-25 14:00:00"
$ws.Range("E4").Value = 7

# Column E was nudged a touch narrower (19 -> 18 characters) in the saved
# workbook; reproduce that via ColumnWidth (stored width uses Excel's
# character/MDW padding, so 17.14 round-trips to a stored width of 18).
$ws.Columns.Item(5).ColumnWidth = 17.14

# The author's cursor ended up on M26 when the workbook was saved.
$ws.Range("M26").Select() | Out-Null
